$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '24.601.09'
Set-TextValue 'E2' '  +3.38%  '
Set-TextValue 'D3' '1.694.90'
Set-TextValue 'E3' '  +1.88%  '
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  +0.28%  '
Set-TextValue 'D5' '316.32'
Set-TextValue 'E5' '  +1.97%  '
Set-TextValue 'E6' '  +0.08%  '
Set-TextValue 'D7' '0.3937'
Set-TextValue 'E7' '  +1.26%  '
Set-TextValue 'D8' '0.4018'
Set-TextValue 'E8' '  +1.81%  '
Set-TextValue 'D9' '1.534'
Set-TextValue 'E9' '  +6.96%  '
Set-TextValue 'B10' 'BinanceUSD'
Set-TextValue 'C10' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D10' '1.001'
Set-TextValue 'E10' '  +0.31%  '
Set-TextValue 'B11' 'OKB'
Set-TextValue 'C11' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D11' '53.89'
Set-TextValue 'E11' '  +9.41%  '
Set-TextValue 'D12' '0.08764'
Set-TextValue 'E12' '  +1.25%  '
Set-TextValue 'D13' '7.217'
Set-TextValue 'E13' '  +7.96%  '
Set-TextValue 'D14' '23.24'
Set-TextValue 'E14' '  +2.48%  '
Set-TextValue 'D15' '0.00001321'
Set-TextValue 'E15' '  +0.50%  '
Set-TextValue 'D16' '7.608'
Set-TextValue 'E16' '  +5.00%  '
Set-TextValue 'D17' '1.699.05'
Set-TextValue 'E17' '  +2.24%  '
Set-TextValue 'D18' '100.18'
Set-TextValue 'E18' '  +0.48%  '
Set-TextValue 'D19' '0.07056'
Set-TextValue 'E19' '  +3.69%  '
Set-TextValue 'D20' '19.65'
Set-TextValue 'E20' '  +3.05%  '
Set-TextValue 'D21' '6.859'
Set-TextValue 'E21' '  +2.84%  '
Set-TextValue 'D22' '1.000'
Set-TextValue 'E22' '  -0.04%  '
Set-TextValue 'D23' '14.04'
Set-TextValue 'E23' '  +1.13%  '
Set-TextValue 'D24' '24.600.68'
Set-TextValue 'E24' '  +3.49%  '
Set-TextValue 'D25' '3.006'
Set-TextValue 'E25' '  +7.31%  '
Set-TextValue 'E26' '  -0.59%  '
Set-TextValue 'D27' '22.37'
Set-TextValue 'E27' '  +2.72%  '
Set-TextValue 'D28' '159.28'
Set-TextValue 'E28' '  +0.54%  '
Set-TextValue 'D29' '5.206'
Set-TextValue 'E29' '  +1.10%  '
Set-TextValue 'D30' '134.05'
Set-TextValue 'E30' '  +3.53%  '
Set-TextValue 'D31' '7.543'
Set-TextValue 'E31' '  +15.94%  '
Set-TextValue 'D32' '1.887.12'
Set-TextValue 'E32' '  +2.38%  '
Set-TextValue 'D33' '1.096'
Set-TextValue 'E33' '  -3.22%  '
Set-TextValue 'D34' '7.315'
Set-TextValue 'E34' '  +11.40%  '
Set-TextValue 'D35' '0.08527'
Set-TextValue 'E35' '  -0.16%  '
Set-TextValue 'D36' '11.36'
Set-TextValue 'E36' '  +9.34%  '
Set-TextValue 'D37' '1.964'
Set-TextValue 'E37' '  +2.14%  '
Set-TextValue 'D38' '0.2723'
Set-TextValue 'E38' '  +2.97%  '
Set-TextValue 'D39' '14.56'
Set-TextValue 'E39' '  +0.57%  '
Set-TextValue 'D40' '0.02760'
Set-TextValue 'E40' '  +9.53%  '
Set-TextValue 'D41' '0.09031'
Set-TextValue 'D42' '1.465'
Set-TextValue 'E42' '  +1.10%  '
Set-TextValue 'D43' '0.7684'
Set-TextValue 'E43' '  +1.67%  '
Set-TextValue 'D44' '0.7183'
Set-TextValue 'E44' '  +2.26%  '
Set-TextValue 'D45' '15.32'
Set-TextValue 'E45' '  +2.77%  '
Set-TextValue 'D46' '2.509'
Set-TextValue 'E46' '  +4.64%  '
Set-TextValue 'D47' '4.205'
Set-TextValue 'E47' '  +2.67%  '
Set-TextValue 'B48' 'Flow'
Set-TextValue 'C48' 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
Set-TextValue 'D48' '1.353'
Set-TextValue 'E48' '  +13.52%  '
Set-TextValue 'B49' 'Frax'
Set-TextValue 'C49' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D49' '1.000'
Set-TextValue 'E49' '  +0.11%  '
Set-TextValue 'D50' '141.39'
Set-TextValue 'E50' '  +2.18%  '
Set-TextValue 'D51' '0.08027'
Set-TextValue 'E51' '  +3.11%  '
